# Adds 2022-Q3 data to the workbook:
#  - Inserts a new row 2 into "总计" (summary) sheet with the 2022-Q3 totals,
#    pushing the existing quarters down by one row.
#  - Inserts a brand-new worksheet named "2022-Q3" (placed before the
#    existing "2022-Q2" sheet) and fills it with the per-fund detail table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert the new 2022-Q3 summary row at row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# Pull the row-index-column style (bold / centered / bordered, same as
# the rest of column A) onto the freshly inserted A2 cell.
$total.Cells.Item(3,1).Copy($total.Cells.Item(2,1))

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 26
$total.Cells.Item(2,4).Value = 2.94

# Clear the borrowed border/bold style from the inserted row's other
# cells (Insert() copies the row-below's formatting by default) so the
# new data cells fall back to the plain/default style like their peers.
$total.Range("B2:D2").ClearFormats()

# Renumber the row-index column (A) for the rows that shifted down so it
# stays a simple 0,1,2,3,4 sequence.
for ($r = 3; $r -le 6; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q3" worksheet with the per-fund holdings detail.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

for ($col = 2; $col -le 8; $col++) {
    # Borrow the bold/centered/bordered header style from the summary
    # sheet's own header row.
    $total.Cells.Item(1,2).Copy($q3.Cells.Item(1,$col))
    $q3.Cells.Item(1,$col).Value = $headers[$col - 2]
}

$rows = @(
    @(0, "501079", "大成科创主题混合（LOF）A", "10.00", "85.01", "6.34", "0.6340", 6),
    @(1, "012473", "大成成长回报六个月持有混合A", "7.48", "80.86", "6.39", "0.4780", 5),
    @(2, "506003", "富国科创板两年定期开放混合", "15.49", "91.82", "2.65", "0.4105", 9),
    @(3, "010371", "大成成长进取混合A", "3.61", "80.75", "6.96", "0.2513", 4),
    @(4, "160642", "鹏华增瑞灵活配置混合（LOF）", "2.05", "90.96", "8.90", "0.1824", 2),
    @(5, "630010", "华商价值精选混合", "4.30", "81.81", "3.92", "0.1686", 4),
    @(6, "011371", "华商远见价值混合型证券投资基金A", "3.85", "64.35", "3.81", "0.1467", 10),
    @(7, "001449", "华商双驱优选灵活配置混合", "2.41", "73.87", "6.08", "0.1465", 1),
    @(8, "010372", "大成成长进取混合C", "1.52", "80.75", "6.96", "0.1058", 4),
    @(9, "020015", "国泰区位优势混合A", "1.95", "87.12", "3.92", "0.0764", 10),
    @(10, "506009", "国泰科创板两年定期开放混合", "2.05", "85.80", "3.49", "0.0715", 10),
    @(11, "002289", "华商改革创新股票A", "1.08", "79.89", "4.37", "0.0472", 2),
    @(12, "630006", "华商产业升级混合", "0.85", "81.97", "3.91", "0.0332", 4),
    @(13, "011073", "鹏华安润混合A", "1.48", "25.50", "2.15", "0.0318", 2),
    @(14, "001723", "华商新动力灵活配置混合", "0.70", "84.40", "3.61", "0.0253", 8),
    @(15, "010403", "华商景气优选混合", "0.61", "76.85", "4.13", "0.0252", 7),
    @(16, "012474", "大成成长回报六个月持有混合C", "0.37", "80.86", "6.39", "0.0236", 5),
    @(17, "004223", "金信多策略精选灵活配置混合", "0.32", "92.79", "6.89", "0.0220", 3),
    @(18, "620002", "金元顺安成长动力混合", "0.51", "62.21", "4.14", "0.0211", 1),
    @(19, "016052", "华商改革创新股票C", "0.32", "79.89", "4.37", "0.0140", 2),
    @(20, "011372", "华商远见价值混合型证券投资基金C", "0.28", "64.35", "3.81", "0.0107", 10),
    @(21, "011886", "弘毅远方高端制造混合型发起式证券投资基金A", "0.25", "89.58", "3.64", "0.0091", 7),
    @(22, "011074", "鹏华安润混合C", "0.19", "25.50", "2.15", "0.0041", 2),
    @(23, "011887", "弘毅远方高端制造混合型发起式证券投资基金C", "0.11", "89.58", "3.64", "0.0040", 7),
    @(24, "016198", "大成科创主题混合（LOF）C", "0.01", "85.01", "6.34", "0.0006", 6),
    @(25, "015594", "国泰区位优势混合C", "0.00", "87.12", "3.92", $null, 10)
)

$r = 2
foreach ($row in $rows) {
    # Column A: plain integer index, styled like the summary sheet's
    # row-index column.
    $total.Cells.Item(3,1).Copy($q3.Cells.Item($r,1))
    $q3.Cells.Item($r,1).Value = $row[0]

    # Columns B-G are stored as text in the source data (fund codes keep
    # their leading zeros, and the numeric-looking figures keep trailing
    # zeros / fixed decimal places) - force text entry via a leading
    # apostrophe, then strip the "stored as text" style back to Normal.
    for ($col = 2; $col -le 6; $col++) {
        $val = $row[$col - 1]
        $q3.Cells.Item($r,$col).Value = "'" + $val
        $q3.Cells.Item($r,$col).Style = "Normal"
    }

    $g = $row[6]
    if ($g -eq $null) {
        # Last row's "持有市值" is a genuine 0 (numeric), not "0.0000" text.
        $q3.Cells.Item($r,7).Value = 0
    } else {
        $q3.Cells.Item($r,7).Value = "'" + $g
        $q3.Cells.Item($r,7).Style = "Normal"
    }

    # Column H: plain numeric rank.
    $q3.Cells.Item($r,8).Value = $row[7]

    $r = $r + 1
}
